$wb = $excel.ActiveWorkbook

# 1) Rename "Categorias_Agua_Potable" to "RAW_Agua_Potable"
$rawWs = $wb.Worksheets.Item("Categorias_Agua_Potable")
$rawWs.Name = "RAW_Agua_Potable"

# 2) Append 5 new rows (431-435) to Categorias_CALL_CENTER
$ccWs = $wb.Worksheets.Item("Categorias_CALL_CENTER")
$ccWs.Cells.Item(431, 1).Value = "CALL CENTER"
$ccWs.Cells.Item(431, 2).Value = "CALL CENTER"
$ccWs.Cells.Item(431, 3).Value = "CATEGORIA 3: OPERACION A 21hs"
$ccWs.Cells.Item(431, 4).Value = "2025-12"
$ccWs.Cells.Item(431, 5).Value = 468939
$ccWs.Cells.Item(431, 6).Value = 17500
$ccWs.Cells.Item(431, 7).Value = 26250
$ccWs.Cells.Item(432, 1).Value = "CALL CENTER"
$ccWs.Cells.Item(432, 2).Value = "CALL CENTER"
$ccWs.Cells.Item(432, 3).Value = "CATEGORIA 3: OPERACION A 21hs"
$ccWs.Cells.Item(432, 4).Value = "2026-01"
$ccWs.Cells.Item(432, 5).Value = 468939
$ccWs.Cells.Item(432, 6).Value = 17500
$ccWs.Cells.Item(432, 7).Value = 26250
$ccWs.Cells.Item(433, 1).Value = "CALL CENTER"
$ccWs.Cells.Item(433, 2).Value = "CALL CENTER"
$ccWs.Cells.Item(433, 3).Value = "CATEGORIA 3: OPERACION A 21hs"
$ccWs.Cells.Item(433, 4).Value = "2026-02"
$ccWs.Cells.Item(433, 5).Value = 468939
$ccWs.Cells.Item(433, 6).Value = 17500
$ccWs.Cells.Item(433, 7).Value = 26250
$ccWs.Cells.Item(434, 1).Value = "CALL CENTER"
$ccWs.Cells.Item(434, 2).Value = "CALL CENTER"
$ccWs.Cells.Item(434, 3).Value = "CATEGORIA 3: OPERACION A 21hs"
$ccWs.Cells.Item(434, 4).Value = "2026-03"
$ccWs.Cells.Item(434, 5).Value = 468939
$ccWs.Cells.Item(434, 6).Value = 17500
$ccWs.Cells.Item(434, 7).Value = 26250
$ccWs.Cells.Item(435, 1).Value = "CALL CENTER"
$ccWs.Cells.Item(435, 2).Value = "CALL CENTER"
$ccWs.Cells.Item(435, 3).Value = "CATEGORIA 3: OPERACION A 21hs"
$ccWs.Cells.Item(435, 4).Value = "2026-04"
$ccWs.Cells.Item(435, 5).Value = 512689
$ccWs.Cells.Item(435, 6).Value = 0
$ccWs.Cells.Item(435, 7).Value = 0

# 3) Add new worksheet "Categorias_AGUA_POTABLE" as the last sheet
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newWs = $wb.Worksheets.Add($null, $lastSheet)
$newCount = $wb.Worksheets.Count
$aguaWs = $wb.Worksheets.Item($newCount)
$aguaWs.Name = "Categorias_AGUA_POTABLE"

# 4) Populate the new sheet with the flattened table data
$aguaWs.Cells.Item(1, 1).Value = "Rama"
$aguaWs.Cells.Item(1, 2).Value = "Agrupamiento"
$aguaWs.Cells.Item(1, 3).Value = "Categoria"
$aguaWs.Cells.Item(1, 4).Value = "Mes"
$aguaWs.Cells.Item(1, 5).Value = "Basico"
$aguaWs.Cells.Item(1, 6).Value = "No_rem"
$aguaWs.Cells.Item(1, 7).Value = "Suma_fija"
$aguaWs.Cells.Item(2, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(2, 2).Value = "PERSONAL SUPERVISIÓN y JEFATURA"
$aguaWs.Cells.Item(2, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(2, 4).Value = "2025-12"
$aguaWs.Cells.Item(2, 5).Value = 3208680
$aguaWs.Cells.Item(2, 6).Value = 120000
$aguaWs.Cells.Item(2, 7).Value = 180000
$aguaWs.Cells.Item(3, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(3, 2).Value = "PERSONAL SUPERVISIÓN y JEFATURA"
$aguaWs.Cells.Item(3, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(3, 4).Value = "2026-01"
$aguaWs.Cells.Item(3, 5).Value = 3208680
$aguaWs.Cells.Item(3, 6).Value = 120000
$aguaWs.Cells.Item(3, 7).Value = 180000
$aguaWs.Cells.Item(4, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(4, 2).Value = "PERSONAL SUPERVISIÓN y JEFATURA"
$aguaWs.Cells.Item(4, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(4, 4).Value = "2026-02"
$aguaWs.Cells.Item(4, 5).Value = 3208680
$aguaWs.Cells.Item(4, 6).Value = 120000
$aguaWs.Cells.Item(4, 7).Value = 180000
$aguaWs.Cells.Item(5, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(5, 2).Value = "PERSONAL SUPERVISIÓN y JEFATURA"
$aguaWs.Cells.Item(5, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(5, 4).Value = "2026-03"
$aguaWs.Cells.Item(5, 5).Value = 3208680
$aguaWs.Cells.Item(5, 6).Value = 120000
$aguaWs.Cells.Item(5, 7).Value = 180000
$aguaWs.Cells.Item(6, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(6, 2).Value = "PERSONAL SUPERVISIÓN y JEFATURA"
$aguaWs.Cells.Item(6, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(6, 4).Value = "2026-04"
$aguaWs.Cells.Item(6, 5).Value = 3508680
$aguaWs.Cells.Item(6, 6).Value = 0
$aguaWs.Cells.Item(6, 7).Value = 0
$aguaWs.Cells.Item(7, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(7, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(7, 3).Value = "OPERADOR DE 2da."
$aguaWs.Cells.Item(7, 4).Value = "2025-12"
$aguaWs.Cells.Item(7, 5).Value = 1764774
$aguaWs.Cells.Item(7, 6).Value = 66000
$aguaWs.Cells.Item(7, 7).Value = 99000
$aguaWs.Cells.Item(8, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(8, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(8, 3).Value = "OPERADOR DE 2da."
$aguaWs.Cells.Item(8, 4).Value = "2026-01"
$aguaWs.Cells.Item(8, 5).Value = 1764774
$aguaWs.Cells.Item(8, 6).Value = 66000
$aguaWs.Cells.Item(8, 7).Value = 99000
$aguaWs.Cells.Item(9, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(9, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(9, 3).Value = "OPERADOR DE 2da."
$aguaWs.Cells.Item(9, 4).Value = "2026-02"
$aguaWs.Cells.Item(9, 5).Value = 1764774
$aguaWs.Cells.Item(9, 6).Value = 66000
$aguaWs.Cells.Item(9, 7).Value = 99000
$aguaWs.Cells.Item(10, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(10, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(10, 3).Value = "OPERADOR DE 2da."
$aguaWs.Cells.Item(10, 4).Value = "2026-03"
$aguaWs.Cells.Item(10, 5).Value = 1764774
$aguaWs.Cells.Item(10, 6).Value = 66000
$aguaWs.Cells.Item(10, 7).Value = 99000
$aguaWs.Cells.Item(11, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(11, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(11, 3).Value = "OPERADOR DE 2da."
$aguaWs.Cells.Item(11, 4).Value = "2026-04"
$aguaWs.Cells.Item(11, 5).Value = 1929774
$aguaWs.Cells.Item(11, 6).Value = 0
$aguaWs.Cells.Item(11, 7).Value = 0
$aguaWs.Cells.Item(12, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(12, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(12, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(12, 4).Value = "2025-12"
$aguaWs.Cells.Item(12, 5).Value = 2032164
$aguaWs.Cells.Item(12, 6).Value = 76000
$aguaWs.Cells.Item(12, 7).Value = 114000
$aguaWs.Cells.Item(13, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(13, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(13, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(13, 4).Value = "2026-01"
$aguaWs.Cells.Item(13, 5).Value = 2032164
$aguaWs.Cells.Item(13, 6).Value = 76000
$aguaWs.Cells.Item(13, 7).Value = 114000
$aguaWs.Cells.Item(14, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(14, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(14, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(14, 4).Value = "2026-02"
$aguaWs.Cells.Item(14, 5).Value = 2032164
$aguaWs.Cells.Item(14, 6).Value = 76000
$aguaWs.Cells.Item(14, 7).Value = 114000
$aguaWs.Cells.Item(15, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(15, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(15, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(15, 4).Value = "2026-03"
$aguaWs.Cells.Item(15, 5).Value = 2032164
$aguaWs.Cells.Item(15, 6).Value = 76000
$aguaWs.Cells.Item(15, 7).Value = 114000
$aguaWs.Cells.Item(16, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(16, 2).Value = "PERSONAL TÉCNICO"
$aguaWs.Cells.Item(16, 3).Value = "OPERADOR DE 1ra."
$aguaWs.Cells.Item(16, 4).Value = "2026-04"
$aguaWs.Cells.Item(16, 5).Value = 2222164
$aguaWs.Cells.Item(16, 6).Value = 0
$aguaWs.Cells.Item(16, 7).Value = 0
$aguaWs.Cells.Item(17, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(17, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(17, 3).Value = "AYUDANTE"
$aguaWs.Cells.Item(17, 4).Value = "2025-12"
$aguaWs.Cells.Item(17, 5).Value = 1176516
$aguaWs.Cells.Item(17, 6).Value = 44000
$aguaWs.Cells.Item(17, 7).Value = 66000
$aguaWs.Cells.Item(18, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(18, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(18, 3).Value = "AYUDANTE"
$aguaWs.Cells.Item(18, 4).Value = "2026-01"
$aguaWs.Cells.Item(18, 5).Value = 1176516
$aguaWs.Cells.Item(18, 6).Value = 44000
$aguaWs.Cells.Item(18, 7).Value = 66000
$aguaWs.Cells.Item(19, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(19, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(19, 3).Value = "AYUDANTE"
$aguaWs.Cells.Item(19, 4).Value = "2026-02"
$aguaWs.Cells.Item(19, 5).Value = 1176516
$aguaWs.Cells.Item(19, 6).Value = 44000
$aguaWs.Cells.Item(19, 7).Value = 66000
$aguaWs.Cells.Item(20, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(20, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(20, 3).Value = "AYUDANTE"
$aguaWs.Cells.Item(20, 4).Value = "2026-03"
$aguaWs.Cells.Item(20, 5).Value = 1176516
$aguaWs.Cells.Item(20, 6).Value = 44000
$aguaWs.Cells.Item(20, 7).Value = 66000
$aguaWs.Cells.Item(21, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(21, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(21, 3).Value = "AYUDANTE"
$aguaWs.Cells.Item(21, 4).Value = "2026-04"
$aguaWs.Cells.Item(21, 5).Value = 1286516
$aguaWs.Cells.Item(21, 6).Value = 0
$aguaWs.Cells.Item(21, 7).Value = 0
$aguaWs.Cells.Item(22, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(22, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(22, 3).Value = "MEDIO OFICIAL / ADMINISTRATIVO 2da."
$aguaWs.Cells.Item(22, 4).Value = "2025-12"
$aguaWs.Cells.Item(22, 5).Value = 1443906
$aguaWs.Cells.Item(22, 6).Value = 54000
$aguaWs.Cells.Item(22, 7).Value = 81000
$aguaWs.Cells.Item(23, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(23, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(23, 3).Value = "MEDIO OFICIAL / ADMINISTRATIVO 2da."
$aguaWs.Cells.Item(23, 4).Value = "2026-01"
$aguaWs.Cells.Item(23, 5).Value = 1443906
$aguaWs.Cells.Item(23, 6).Value = 54000
$aguaWs.Cells.Item(23, 7).Value = 81000
$aguaWs.Cells.Item(24, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(24, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(24, 3).Value = "MEDIO OFICIAL / ADMINISTRATIVO 2da."
$aguaWs.Cells.Item(24, 4).Value = "2026-02"
$aguaWs.Cells.Item(24, 5).Value = 1443906
$aguaWs.Cells.Item(24, 6).Value = 54000
$aguaWs.Cells.Item(24, 7).Value = 81000
$aguaWs.Cells.Item(25, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(25, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(25, 3).Value = "MEDIO OFICIAL / ADMINISTRATIVO 2da."
$aguaWs.Cells.Item(25, 4).Value = "2026-03"
$aguaWs.Cells.Item(25, 5).Value = 1443906
$aguaWs.Cells.Item(25, 6).Value = 54000
$aguaWs.Cells.Item(25, 7).Value = 81000
$aguaWs.Cells.Item(26, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(26, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(26, 3).Value = "MEDIO OFICIAL / ADMINISTRATIVO 2da."
$aguaWs.Cells.Item(26, 4).Value = "2026-04"
$aguaWs.Cells.Item(26, 5).Value = 1578906
$aguaWs.Cells.Item(26, 6).Value = 0
$aguaWs.Cells.Item(26, 7).Value = 0
$aguaWs.Cells.Item(27, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(27, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(27, 3).Value = "OFICIAL / ADMINISTRATIVO 1ra."
$aguaWs.Cells.Item(27, 4).Value = "2025-12"
$aguaWs.Cells.Item(27, 5).Value = 1604340
$aguaWs.Cells.Item(27, 6).Value = 60000
$aguaWs.Cells.Item(27, 7).Value = 90000
$aguaWs.Cells.Item(28, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(28, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(28, 3).Value = "OFICIAL / ADMINISTRATIVO 1ra."
$aguaWs.Cells.Item(28, 4).Value = "2026-01"
$aguaWs.Cells.Item(28, 5).Value = 1604340
$aguaWs.Cells.Item(28, 6).Value = 60000
$aguaWs.Cells.Item(28, 7).Value = 90000
$aguaWs.Cells.Item(29, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(29, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(29, 3).Value = "OFICIAL / ADMINISTRATIVO 1ra."
$aguaWs.Cells.Item(29, 4).Value = "2026-02"
$aguaWs.Cells.Item(29, 5).Value = 1604340
$aguaWs.Cells.Item(29, 6).Value = 60000
$aguaWs.Cells.Item(29, 7).Value = 90000
$aguaWs.Cells.Item(30, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(30, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(30, 3).Value = "OFICIAL / ADMINISTRATIVO 1ra."
$aguaWs.Cells.Item(30, 4).Value = "2026-03"
$aguaWs.Cells.Item(30, 5).Value = 1604340
$aguaWs.Cells.Item(30, 6).Value = 60000
$aguaWs.Cells.Item(30, 7).Value = 90000
$aguaWs.Cells.Item(31, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(31, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(31, 3).Value = "OFICIAL / ADMINISTRATIVO 1ra."
$aguaWs.Cells.Item(31, 4).Value = "2026-04"
$aguaWs.Cells.Item(31, 5).Value = 1754340
$aguaWs.Cells.Item(31, 6).Value = 0
$aguaWs.Cells.Item(31, 7).Value = 0
$aguaWs.Cells.Item(32, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(32, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(32, 3).Value = "OFICIAL ENCARGADO / ENCARGADO"
$aguaWs.Cells.Item(32, 4).Value = "2025-12"
$aguaWs.Cells.Item(32, 5).Value = 1818252
$aguaWs.Cells.Item(32, 6).Value = 68000
$aguaWs.Cells.Item(32, 7).Value = 102000
$aguaWs.Cells.Item(33, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(33, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(33, 3).Value = "OFICIAL ENCARGADO / ENCARGADO"
$aguaWs.Cells.Item(33, 4).Value = "2026-01"
$aguaWs.Cells.Item(33, 5).Value = 1818252
$aguaWs.Cells.Item(33, 6).Value = 68000
$aguaWs.Cells.Item(33, 7).Value = 102000
$aguaWs.Cells.Item(34, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(34, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(34, 3).Value = "OFICIAL ENCARGADO / ENCARGADO"
$aguaWs.Cells.Item(34, 4).Value = "2026-02"
$aguaWs.Cells.Item(34, 5).Value = 1818252
$aguaWs.Cells.Item(34, 6).Value = 68000
$aguaWs.Cells.Item(34, 7).Value = 102000
$aguaWs.Cells.Item(35, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(35, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(35, 3).Value = "OFICIAL ENCARGADO / ENCARGADO"
$aguaWs.Cells.Item(35, 4).Value = "2026-03"
$aguaWs.Cells.Item(35, 5).Value = 1818252
$aguaWs.Cells.Item(35, 6).Value = 68000
$aguaWs.Cells.Item(35, 7).Value = 102000
$aguaWs.Cells.Item(36, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(36, 2).Value = "PERSONAL AUXILIAR / ADMINISTRATIVO"
$aguaWs.Cells.Item(36, 3).Value = "OFICIAL ENCARGADO / ENCARGADO"
$aguaWs.Cells.Item(36, 4).Value = "2026-04"
$aguaWs.Cells.Item(36, 5).Value = 1988252
$aguaWs.Cells.Item(36, 6).Value = 0
$aguaWs.Cells.Item(36, 7).Value = 0
$aguaWs.Cells.Item(37, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(37, 2).Value = "MAESTRANZA"
$aguaWs.Cells.Item(37, 3).Value = "Maestranza C"
$aguaWs.Cells.Item(37, 4).Value = "2025-12"
$aguaWs.Cells.Item(37, 5).Value = 1069560
$aguaWs.Cells.Item(37, 6).Value = 40000
$aguaWs.Cells.Item(37, 7).Value = 60000
$aguaWs.Cells.Item(38, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(38, 2).Value = "MAESTRANZA"
$aguaWs.Cells.Item(38, 3).Value = "Maestranza C"
$aguaWs.Cells.Item(38, 4).Value = "2026-01"
$aguaWs.Cells.Item(38, 5).Value = 1069560
$aguaWs.Cells.Item(38, 6).Value = 40000
$aguaWs.Cells.Item(38, 7).Value = 60000
$aguaWs.Cells.Item(39, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(39, 2).Value = "MAESTRANZA"
$aguaWs.Cells.Item(39, 3).Value = "Maestranza C"
$aguaWs.Cells.Item(39, 4).Value = "2026-02"
$aguaWs.Cells.Item(39, 5).Value = 1069560
$aguaWs.Cells.Item(39, 6).Value = 40000
$aguaWs.Cells.Item(39, 7).Value = 60000
$aguaWs.Cells.Item(40, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(40, 2).Value = "MAESTRANZA"
$aguaWs.Cells.Item(40, 3).Value = "Maestranza C"
$aguaWs.Cells.Item(40, 4).Value = "2026-03"
$aguaWs.Cells.Item(40, 5).Value = 1069560
$aguaWs.Cells.Item(40, 6).Value = 40000
$aguaWs.Cells.Item(40, 7).Value = 60000
$aguaWs.Cells.Item(41, 1).Value = "AGUA POTABLE"
$aguaWs.Cells.Item(41, 2).Value = "MAESTRANZA"
$aguaWs.Cells.Item(41, 3).Value = "Maestranza C"
$aguaWs.Cells.Item(41, 4).Value = "2026-04"
$aguaWs.Cells.Item(41, 5).Value = 1169560
$aguaWs.Cells.Item(41, 6).Value = 0
$aguaWs.Cells.Item(41, 7).Value = 0

# 5) Restore original active sheet/selection (Indice) since Worksheets.Add()
#    switches focus to the newly created sheet
$wb.Worksheets.Item(1).Activate()
